$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.303.99'
$ws.Range('E2').Value = '  +2.99%  '
$ws.Range('D3').Value = '2.264.47'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'321.32"
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = "'102.87"
$ws.Range('E6').Value = '  +4.29%  '
$ws.Range('D7').Value = "'0.580"
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = "'0.557"
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').Value = "'37.68"
$ws.Range('E10').Value = '  +2.48%  '
$ws.Range('D11').Value = "'0.0836"
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('D12').Value = "'7.65"
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = "'0.107"
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').Value = '2.599.97'
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = "'14.56"
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = "'0.863"
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').Value = '2.263.07'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').Value = '44.141.71'
$ws.Range('E18').Value = '  +2.92%  '
$ws.Range('D19').Value = "'13.41"
$ws.Range('E19').Value = '  -4.07%  '
$ws.Range('D20').Value = '0.0₃0987'
$ws.Range('E20').Value = '  +2.51%  '
$ws.Range('D21').Value = "'6.55"
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').Value = "'65.91"
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').Value = "'3.17"
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('D24').Value = "'236.72"
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').Value = "'2.13"
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').Value = "'10.29"
$ws.Range('E27').Value = '  +3.26%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = "'38.47"
$ws.Range('E28').Value = '  +7.34%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'2.20"
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('D30').Value = "'6.24"
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('D31').Value = "'162.92"
$ws.Range('E31').Value = '  +6.04%  '
$ws.Range('D32').Value = "'20.32"
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').Value = "'0.0858"
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').Value = "'0.116"
$ws.Range('E35').Value = '  +11.23%  '
$ws.Range('D36').Value = "'1.97"
$ws.Range('E36').Value = '  +2.77%  '
$ws.Range('D37').Value = "'3.05"
$ws.Range('E37').Value = '  -5.01%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').Value = "'17.15"
$ws.Range('E38').Value = '  +24.34%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = "'0.120"
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('D40').Value = "'3.76"
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('D41').Value = "'4.23"
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('D42').Value = "'0.0318"
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').Value = '1.790.93'
$ws.Range('E44').Value = '  +3.95%  '
$ws.Range('B45').Value = 'ordi'
$ws.Range('C45').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D45').Value = "'76.77"
$ws.Range('E45').Value = '  +3.01%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = "'0.200"
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('D47').Value = "'83.26"
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'15.92"
$ws.Range('E48').Value = '  +6.98%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = "'5.25"
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('D50').Value = "'1.71"
$ws.Range('E50').Value = '  +9.16%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'105.37"
$ws.Range('E51').Value = '  +2.54%  '
